# Apply updated dSF (column F) values for the specified rows.
# Mapping of row number -> new value, per the commit's repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -3
    7  = -1
    8  = -6
    10 = -2
    11 = 1
    15 = -3
    16 = 3
    20 = 0
    21 = -1
    22 = 7
    28 = -4
    33 = -3
    34 = -3
    35 = 8
    36 = -2
    42 = -15
    43 = -2
    46 = 1
    51 = 7
    59 = 1
    61 = -3
    63 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
